# Refresh the scraped crypto price/volume snapshot (coinranking.com) to the
# latest GitHub Actions run. Two rows swap rank position in this update:
# WrappedEther now outranks TRON (rows 18-19), and Kaspa now outranks
# TheGraph (rows 43-44); their Coin/Link/Price/Volume cells are rewritten
# together below. All other rows keep their rank but get fresh Price /
# Volume(1h) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v = '66.806.10'
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -5.35%  '
$cell = $ws.Cells.Item(2, 5)
$cell.Value = $v
$v = '3.369.54'
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -6.65%  '
$cell = $ws.Cells.Item(3, 5)
$cell.Value = $v
$v = '  +0.04%  '
$cell = $ws.Cells.Item(4, 5)
$cell.Value = $v
$v = '561.52'
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -5.87%  '
$cell = $ws.Cells.Item(5, 5)
$cell.Value = $v
$v = '183.98'
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -9.54%  '
$cell = $ws.Cells.Item(6, 5)
$cell.Value = $v
$v = '0.596'
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -5.09%  '
$cell = $ws.Cells.Item(7, 5)
$cell.Value = $v
$v = '  +0.02%  '
$cell = $ws.Cells.Item(8, 5)
$cell.Value = $v
$v = '3.364.48'
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -6.44%  '
$cell = $ws.Cells.Item(9, 5)
$cell.Value = $v
$v = '0.189'
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -12.51%  '
$cell = $ws.Cells.Item(10, 5)
$cell.Value = $v
$v = '0.597'
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -7.46%  '
$cell = $ws.Cells.Item(11, 5)
$cell.Value = $v
$v = '48.01'
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -11.07%  '
$cell = $ws.Cells.Item(12, 5)
$cell.Value = $v
$v = '0.0000269'
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -10.99%  '
$cell = $ws.Cells.Item(13, 5)
$cell.Value = $v
$v = '8.76'
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -9.41%  '
$cell = $ws.Cells.Item(14, 5)
$cell.Value = $v
$v = '3.904.53'
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -6.58%  '
$cell = $ws.Cells.Item(15, 5)
$cell.Value = $v
$v = '611.55'
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -10.37%  '
$cell = $ws.Cells.Item(16, 5)
$cell.Value = $v
$v = '66.676.83'
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -5.60%  '
$cell = $ws.Cells.Item(17, 5)
$cell.Value = $v
$v = 'WrappedEther'
$cell = $ws.Cells.Item(18, 2)
$cell.Value = $v
$v = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$cell = $ws.Cells.Item(18, 3)
$cell.Value = $v
$v = '3.375.20'
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -6.62%  '
$cell = $ws.Cells.Item(18, 5)
$cell.Value = $v
$v = 'TRON'
$cell = $ws.Cells.Item(19, 2)
$cell.Value = $v
$v = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$cell = $ws.Cells.Item(19, 3)
$cell.Value = $v
$v = '0.118'
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -3.79%  '
$cell = $ws.Cells.Item(19, 5)
$cell.Value = $v
$v = '17.72'
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -7.50%  '
$cell = $ws.Cells.Item(20, 5)
$cell.Value = $v
$v = '11.73'
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -8.16%  '
$cell = $ws.Cells.Item(21, 5)
$cell.Value = $v
$v = '0.914'
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -8.39%  '
$cell = $ws.Cells.Item(22, 5)
$cell.Value = $v
$v = '17.05'
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -7.53%  '
$cell = $ws.Cells.Item(23, 5)
$cell.Value = $v
$v = '5.10'
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -3.34%  '
$cell = $ws.Cells.Item(24, 5)
$cell.Value = $v
$v = '96.25'
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -12.90%  '
$cell = $ws.Cells.Item(25, 5)
$cell.Value = $v
$v = '4.09'
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -10.13%  '
$cell = $ws.Cells.Item(26, 5)
$cell.Value = $v
$v = '2.76'
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -8.87%  '
$cell = $ws.Cells.Item(27, 5)
$cell.Value = $v
$v = '9.60'
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -9.70%  '
$cell = $ws.Cells.Item(28, 5)
$cell.Value = $v
$v = '8.84'
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -12.38%  '
$cell = $ws.Cells.Item(29, 5)
$cell.Value = $v
$v = '30.92'
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -10.11%  '
$cell = $ws.Cells.Item(30, 5)
$cell.Value = $v
$v = '6.38'
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -11.24%  '
$cell = $ws.Cells.Item(31, 5)
$cell.Value = $v
$v = '3.85'
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -13.73%  '
$cell = $ws.Cells.Item(32, 5)
$cell.Value = $v
$v = '11.27'
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -8.76%  '
$cell = $ws.Cells.Item(33, 5)
$cell.Value = $v
$v = '0.106'
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -7.47%  '
$cell = $ws.Cells.Item(34, 5)
$cell.Value = $v
$v = '58.73'
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -7.69%  '
$cell = $ws.Cells.Item(35, 5)
$cell.Value = $v
$v = '3.795.36'
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -1.71%  '
$cell = $ws.Cells.Item(36, 5)
$cell.Value = $v
$v = '532.96'
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  +4.52%  '
$cell = $ws.Cells.Item(37, 5)
$cell.Value = $v
$v = '1.00'
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  +0.04%  '
$cell = $ws.Cells.Item(38, 5)
$cell.Value = $v
$v = '3.77'
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  +37.37%  '
$cell = $ws.Cells.Item(39, 5)
$cell.Value = $v
$v = '3.44'
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -4.15%  '
$cell = $ws.Cells.Item(40, 5)
$cell.Value = $v
$v = '0.0{0}0730' -f [char]0x2083
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -14.30%  '
$cell = $ws.Cells.Item(41, 5)
$cell.Value = $v
$v = '  -9.48%  '
$cell = $ws.Cells.Item(42, 5)
$cell.Value = $v
$v = 'Kaspa'
$cell = $ws.Cells.Item(43, 2)
$cell.Value = $v
$v = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell = $ws.Cells.Item(43, 3)
$cell.Value = $v
$v = '0.128'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -6.64%  '
$cell = $ws.Cells.Item(43, 5)
$cell.Value = $v
$v = 'TheGraph'
$cell = $ws.Cells.Item(44, 2)
$cell.Value = $v
$v = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$cell = $ws.Cells.Item(44, 3)
$cell.Value = $v
$v = '0.353'
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -8.17%  '
$cell = $ws.Cells.Item(44, 5)
$cell.Value = $v
$v = '32.86'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -10.74%  '
$cell = $ws.Cells.Item(45, 5)
$cell.Value = $v
$v = '0.0419'
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -11.05%  '
$cell = $ws.Cells.Item(46, 5)
$cell.Value = $v
$v = '3.19'
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -7.57%  '
$cell = $ws.Cells.Item(47, 5)
$cell.Value = $v
$v = '  -12.60%  '
$cell = $ws.Cells.Item(48, 5)
$cell.Value = $v
$v = '0.131'
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -7.27%  '
$cell = $ws.Cells.Item(49, 5)
$cell.Value = $v
$v = '0.998'
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -0.32%  '
$cell = $ws.Cells.Item(50, 5)
$cell.Value = $v
$v = '7.78'
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = $v
$cell.Style = "Normal"
$v = '  -10.07%  '
$cell = $ws.Cells.Item(51, 5)
$cell.Value = $v
